# Penalty Reward System (unfinished) - remove the two weekly rows that fell
# within the June 2023 order month (45088.99999999999 / 15 and
# 45095.99999999999 / 14) from the "Weekly Quantity" sheet, shifting all
# rows below them up by two, and correspondingly reduce the June 2023
# total on the "Monthly Trend" sheet from 50 to 21 (50 - 15 - 14).

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Delete worksheet rows 5 and 6 (dates 45088.99999999999/15 and
# 45095.99999999999/14), shifting everything below up by two rows.
$wsWeekly.Range("A5:B6").EntireRow.Delete()

# Update the corresponding monthly total (June 2023) to reflect the
# removed weekly quantities.
$wsMonthly.Range("B4").Value = 21
